# netCrypto.xlsx — "Add files via upload" re-edit:
#   - T2 / T3 USD-amount values bumped from 0 to 10
#   - view scrolled right so column O is the left-most visible column
#   - selection moved to T4

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Make sure we're working on the sheet that is being edited.
$ws.Activate()

# --- data edits -----------------------------------------------------------
$ws.Range("T2").Value = 10
$ws.Range("T3").Value = 10

# --- view state -------------------------------------------------------
# Select T4 and scroll the window so O1 becomes the top-left visible cell.
$ws.Range("T4").Select()
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
